# Updates the "cryptos" price/volume table (columns D = Price, E = Volume(1h))
# to the refreshed values from the Sat Nov 30 21:46:22 UTC 2024 GitHub
# Actions run. Cells that look like a plain decimal number (e.g. "1.93")
# are forced to Text via NumberFormat "@" before the write and then the
# style is reset to "Normal" afterwards, so Excel's COM layer doesn't
# silently convert them to a numeric value (which would change both the
# stored type and the displayed precision) while leaving no stray style
# behind on the cell itself.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.836.67'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '3.701.10'
$ws.Range('E3').Value = '  +3.25%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.84'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +9.32%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '657.08'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.74%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.428'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('E9').Value = '  +3.80%  '
$ws.Range('E10').Value = '  +0.04%  '
$ws.Range('D11').Value = '3.698.31'
$ws.Range('E11').Value = '  +3.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '45.43'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.86%  '
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('E14').Value = '  +6.53%  '
$ws.Range('D15').Value = '4.390.92'
$ws.Range('E16').Value = '  +3.79%  '
$ws.Range('D17').Value = '96.639.28'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.08'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +17.04%  '
$ws.Range('D19').Value = '3.688.49'
$ws.Range('E19').Value = '  +2.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.89'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.535'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '525.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.51'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '7.10'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.31%  '
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '102.72'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '13.45'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.169'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.30%  '
$ws.Range('E30').Value = '  +5.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.08'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.18%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('E33').Value = '  +14.63%  '
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('E35').Value = '  +2.88%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '656.55'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +6.50%  '
$ws.Range('E38').Value = '  +4.40%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.99'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.18'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +17.77%  '
$ws.Range('E41').Value = '  +5.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.01'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.27%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.972'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.85'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +18.88%  '
$ws.Range('E46').Value = '  +3.76%  '
$ws.Range('E47').Value = '  +1.92%  '
$ws.Range('E48').Value = '  +0.84%  '
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('E50').Value = '  +0.43%  '
$ws.Range('E51').Value = '  +2.66%  '
